# Remove the "课件" (teaching materials) section from the resource sheet.
#
# The sheet is organised as a sequence of stacked sections, each consisting
# of a merged title row, a merged column-header row, and a handful of blank
# data rows. The "课件" section occupies rows 30-36 (title row 30, header
# row 31, blank rows 32-36) and is immediately followed by the "教材"
# (textbooks) section in rows 37-43 (title row 37, header row 38, blank
# rows 39-43).
#
# Deleting the "课件" section's rows shifts everything below it (the
# "教材" section) upward by 7 rows, so "教材" ends up occupying rows
# 30-36 - exactly matching the target layout. Excel's own row-delete
# machinery also takes care of: compacting the shared-string table (the
# now-unused "文件夹"/"授课老师"/"文件个数"/"课件" strings are dropped),
# renumbering every remaining cell reference, shrinking the merged-cell
# list (the A37:H37 title-row merge disappears), and narrowing the
# "文件夹命名" data-validation sqref from "B32:B36 B39:B43" down to just
# "B32:B36".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30:H36").EntireRow.Delete() | Out-Null

$ws.Range("J32").Select() | Out-Null
